$wb = $excel.ActiveWorkbook

# pir sheet - rows 224-236
$wsPIR = $wb.Worksheets.Item("PIR")
$pirDates = @("2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28")
$pirTimes = @("16:55:44","16:55:48","16:55:50","16:55:55","16:56:00","16:56:05","16:56:10","16:56:15","16:56:20","16:56:26","16:56:30","16:56:35","16:56:40")
$pirHours = @("16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00")
$pirLocations = @("Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom")
$pirValues = @("No Motion","No Motion","No Motion","No Motion","No Motion","No Motion","No Motion","No Motion","No Motion","No Motion","No Motion","No Motion","No Motion")
$pirStatuses = @("Inactive","Inactive","Inactive","Inactive","Inactive","Inactive","Inactive","Inactive","Inactive","Inactive","Inactive","Inactive","Inactive")
$wsPIR.Range("A224`:A236").NumberFormat = "@"
for ($i = 0; $i -lt $pirDates.Length; $i++) {
    $r = 224 + $i
    $wsPIR.Cells.Item($r, 1).Value = $pirDates[$i]
    $wsPIR.Cells.Item($r, 2).Value = $pirTimes[$i]
    $wsPIR.Cells.Item($r, 3).Value = $pirHours[$i]
    $wsPIR.Cells.Item($r, 4).Value = $pirLocations[$i]
    $wsPIR.Cells.Item($r, 5).Value = $pirValues[$i]
    $wsPIR.Cells.Item($r, 6).Value = $pirStatuses[$i]
}

# hum sheet - rows 218-232
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humDates = @("2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28")
$humTimes = @("16:55:45","16:55:46","16:55:52","16:55:56","16:56:00","16:56:04","16:56:08","16:56:16","16:56:20","16:56:24","16:56:28","16:56:32","16:56:36","16:56:41","16:56:44")
$humHours = @("16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00")
$humLocations = @("Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom")
$humValues = @("87.1%","88.0%","88.0%","88.0%","87.9%","87.9%","87.0%","87.9%","86.9%","87.9%","87.0%","87.9%","87.9%","86.9%","87.9%")
$humStatuses = @("Active","Active","Active","Active","Active","Active","Active","Active","Active","Active","Active","Active","Active","Active","Active")
$wsHumidity.Range("A218`:A232").NumberFormat = "@"
$wsHumidity.Range("E218`:E232").NumberFormat = "@"
for ($i = 0; $i -lt $humDates.Length; $i++) {
    $r = 218 + $i
    $wsHumidity.Cells.Item($r, 1).Value = $humDates[$i]
    $wsHumidity.Cells.Item($r, 2).Value = $humTimes[$i]
    $wsHumidity.Cells.Item($r, 3).Value = $humHours[$i]
    $wsHumidity.Cells.Item($r, 4).Value = $humLocations[$i]
    $wsHumidity.Cells.Item($r, 5).Value = $humValues[$i]
    $wsHumidity.Cells.Item($r, 6).Value = $humStatuses[$i]
}

# temp sheet - rows 218-231
$wsTemperature = $wb.Worksheets.Item("Temperature")
$tempDates = @("2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28","2026-01-28")
$tempTimes = @("16:55:46","16:55:47","16:55:53","16:55:57","16:56:01","16:56:05","16:56:09","16:56:17","16:56:21","16:56:25","16:56:29","16:56:33","16:56:37","16:56:41")
$tempHours = @("16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00","16:00")
$tempLocations = @("Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom","Bathroom")
$tempValues = @("22.8C","22.8C","22.8C","22.8C","22.8C","22.8C","22.8C","22.8C","22.7C","22.8C","22.8C","22.8C","22.8C","22.8C")
$tempStatuses = @("Active","Active","Active","Active","Active","Active","Active","Active","Active","Active","Active","Active","Active","Active")
$wsTemperature.Range("A218`:A231").NumberFormat = "@"
for ($i = 0; $i -lt $tempDates.Length; $i++) {
    $r = 218 + $i
    $wsTemperature.Cells.Item($r, 1).Value = $tempDates[$i]
    $wsTemperature.Cells.Item($r, 2).Value = $tempTimes[$i]
    $wsTemperature.Cells.Item($r, 3).Value = $tempHours[$i]
    $wsTemperature.Cells.Item($r, 4).Value = $tempLocations[$i]
    $wsTemperature.Cells.Item($r, 5).Value = $tempValues[$i]
    $wsTemperature.Cells.Item($r, 6).Value = $tempStatuses[$i]
}
